# "move name to middle initial field"
#
# Personnel!A5 currently holds "E. Taylor" as a single givenName value with
# the middleInitial column (B5) empty. Split it: the surname-like "Taylor"
# part moves into the middleInitial field (B5), leaving the initial "E."
# behind in the givenName field (A5). surName (C5, "Crockford") is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Write the middle-initial field first, then trim the given-name field, so
# the new shared strings are interned in the same order as the source edit.
$ws.Range("B5").Value = "Taylor"
$ws.Range("A5").Value = "E."

# The editor ended up focused on the Personnel sheet with A5 selected.
$ws.Activate() | Out-Null
$ws.Range("A5").Select() | Out-Null
